$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '26.416.25'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.615.57'
$ws.Range('E3').Value = '  +1.66%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '214.00'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').Value = '0.503'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.0610'
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = '19.20'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('D11').Value = '0.0855'
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('D12').Value = '1.839.59'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').Value = '1.627.10'
$ws.Range('E13').Value = '  +4.12%  '
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').Value = '64.69'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').Value = '26.416.56'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').Value = '225.46'
$ws.Range('E18').Value = '  +6.01%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0727'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '7.53'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').Value = '4.37'
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('D23').Value = '9.09'
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = '145.10'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').Value = '7.04'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D28').Value = '0.114'
$ws.Range('E28').Value = '  +1.85%  '
$ws.Range('D29').Value = '15.33'
$ws.Range('E29').Value = '  +0.90%  '
$ws.Range('D30').Value = '0.0498'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('D32').Value = '3.22'
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').Value = '1.448.38'
$ws.Range('E33').Value = '  +8.57%  '
$ws.Range('D34').Value = '2.99'
$ws.Range('E34').Value = '  +1.84%  '
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('E37').Value = '  -4.92%  '
$ws.Range('D38').Value = '0.0167'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('D39').Value = '0.837'
$ws.Range('E39').Value = '  +2.24%  '
$ws.Range('E40').Value = '  +2.03%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = '2.19'
$ws.Range('E42').Value = '  +2.33%  '
$ws.Range('D43').Value = '1.751.45'
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('D44').Value = '0.761'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '61.90'
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').Value = '0.911'
$ws.Range('E46').Value = '  -10.81%  '
$ws.Range('D47').Value = '88.01'
$ws.Range('E47').Value = '  +2.80%  '
$ws.Range('D48').Value = '0.0₆0109'
$ws.Range('E48').Value = '  +4.82%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  -1.16%  '
